$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (e.g. "592.94"), so they stay text
# like the other inline-string price cells (e.g. "65.111.31").
$textCells = @("D5","D6","D14","D19","D21","D22","D25","D27","D28","D29","D30","D31","D32","D37","D39","D40","D42","D43","D45","D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "65.111.31"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "3.523.59"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "592.94"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("D6").Value = "134.32"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("D7").Value = "3.522.34"
$ws.Range("E7").Value = "  -1.28%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("E11").Value = "  +2.68%  "
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "4.122.67"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").Value = "27.63"
$ws.Range("E14").Value = "  +1.48%  "
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").Value = "3.531.56"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "65.047.81"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").Value = "10.06"
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").Value = "5.68"
$ws.Range("E21").Value = "  -3.30%  "
$ws.Range("D22").Value = "392.59"
$ws.Range("E22").Value = "  +1.04%  "
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").Value = "3.667.21"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("D25").Value = "74.75"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "0.0000112"
$ws.Range("E27").Value = "  -4.59%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "7.69"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").Value = "1.59"
$ws.Range("E29").Value = "  +7.75%  "
$ws.Range("D30").Value = "0.996"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").Value = "2.27"
$ws.Range("E31").Value = "  -1.74%  "
$ws.Range("D32").Value = "8.35"
$ws.Range("E32").Value = "  -1.40%  "
$ws.Range("D33").Value = "3.532.32"
$ws.Range("E33").Value = "  -1.20%  "
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").Value = "5.31"
$ws.Range("E37").Value = "  +5.18%  "
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("D39").Value = "6.94"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("D40").Value = "168.41"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").Value = "0.823"
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("D43").Value = "1.26"
$ws.Range("E43").Value = "  +4.44%  "
$ws.Range("E44").Value = "  +0.73%  "
$ws.Range("D45").Value = "25.62"
$ws.Range("E45").Value = "  -6.36%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("D50").Value = "2.421.56"
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("D51").Value = "0.903"
$ws.Range("E51").Value = "  +4.25%  "

# Restore default (General) style on the cells we forced to text so we
# do not leave a residual per-cell style reference.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
